$d = $word.ActiveDocument

function Remove-CharAfter($searchText) {
    # Finds the first occurrence of $searchText (from the start of the
    # document) and deletes the single character immediately following the
    # match (used to drop a trailing, standalone "," or "." run).
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }
    $tail = $d.Range($r.End, $r.End + 1)
    $tail.Delete()
}

# 1. Add the team members' names after "Mandataire : DeGuiWii"
$d.Content.Find.Execute(" DeGuiWii", $true, $false, $false, $false, $false, $true, 1, $false, " DeGuiWii (Denis Thériault, Guillaume Gagnon, William Lemieux)", 2) | Out-Null

# 2. Mention "des employés de terrain" before the "responsable d'un projet" clause
$d.Content.Find.Execute(", responsable d’un projet, entre ses heures sur le site", $true, $false, $false, $false, $false, $true, 1, $false, " des employés de terrain, responsable d’un projet, entre ses heures sur le site", 2) | Out-Null

# 3-17. Drop the trailing comma/period runs that turned the bullet list into
# a plain list without closing punctuation.
Remove-CharAfter("Gestion des utilisateurs")
Remove-CharAfter("Gestion des projets")
Remove-CharAfter("Gestion des feuilles de temps")
Remove-CharAfter(" courriels")
Remove-CharAfter("Création d’Excel pour des aperçus (rapports, graphiques)")
Remove-CharAfter("Module de gestion de feuille de temps")
Remove-CharAfter(" et graphiques)")
Remove-CharAfter("de gestion des employés")
Remove-CharAfter("Module des calculs")
Remove-CharAfter("Module de recherche de projet")
Remove-CharAfter("Module compte de dépense")
Remove-CharAfter("données")
Remove-CharAfter("Calculer le temps investi dans chaque projet")
Remove-CharAfter("Calculer les heures travaillées")
Remove-CharAfter("es dépenses d’un projet")

# 18. "les lundis matin" -> "les dimanches soir"
$d.Content.Find.Execute("nvoi de courriel les lundis matin aux employés", $true, $false, $false, $false, $false, $true, 1, $false, "nvoi de courriel les dimanches soir aux employés", 2) | Out-Null

# 19. Capitalize "avoir" -> "Avoir"
$d.Content.Find.Execute("(avoir un historique des projets)", $true, $false, $false, $false, $false, $true, 1, $false, "(Avoir un historique des projets)", 2) | Out-Null
